# Rework the "name/province/tazkira" table into a "country/year" table
# (mirrors the excel.py rewrite from pandas -> polars: the demographic
# sample data was swapped for a small country/year list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old third column (tazkira / numeric id) entirely - the new
# table only has two columns.
$ws.Range("C1:C2").Delete()

# New header row.
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "year"

# New data rows.
$ws.Range("A2").Value = "Afghanistan"
$ws.Range("B2").Value = 2025
$ws.Range("A3").Value = "Pakistan"
$ws.Range("B3").Value = 2025
$ws.Range("A4").Value = "Iran"
$ws.Range("B4").Value = 2025

# Keep the new rows' formatting consistent with the existing data row
# (same cell style as row 2) instead of picking up a blank default style.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B4").PasteSpecial(-4122)
